# Updates for 2 Apr
# Refresh the COVID-19 "deaths" figures (column C, Sheet1) with the
# latest counts reported on 2 April.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row (on Sheet1) -> new deaths value
$updates = @{
    2  = 28    # Alabama
    4  = 29    # Arizona
    5  = 10    # Arkansas
    6  = 215   # California
    7  = 80    # Colorado
    8  = 85    # Connecticut
    9  = 11    # Delaware
    10 = 11    # District of Columbia
    11 = 101   # Florida
    12 = 154   # Georgia
    16 = 141   # Illinois
    17 = 65    # Indiana
    18 = 9     # Iowa
    19 = 11    # Kansas
    20 = 20    # Kentucky
    21 = 273   # Louisiana
    22 = 7     # Maine
    23 = 31    # Maryland
    24 = 122   # Massachusetts
    25 = 337   # Michigan
    26 = 17    # Minnesota
    27 = 22    # Mississippi
    28 = 18    # Missouri
    29 = 6     # Montana
    30 = 4     # Nebraska
    31 = 26    # Nevada
    32 = 4     # New Hampshire
    33 = 355   # New Jersey
    34 = 6     # New Mexico
    35 = 2219  # New York
    36 = 15    # North Carolina
    38 = 65    # Ohio
    39 = 30    # Oklahoma
    40 = 19    # Oregon
    41 = 74    # Pennsylvania
    42 = 12    # Puerto Rico
    43 = 10    # Rhode Island
    44 = 26    # South Carolina
    45 = 2     # South Dakota
    46 = 25    # Tennessee
    47 = 60    # Texas
    48 = 7     # Utah
    49 = 16    # Vermont
    50 = 34    # Virginia
    51 = 254   # Washington
    52 = 2     # West Virginia
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Reproduce the saved view state as closely as possible: scrolled back to
# the top (no frozen/offset topLeftCell) with C2:C54 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C2:C54").Select()
